$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lexiconv2_lematizado")

# Append new lexicon rows (426-529) with words + binary flags in columns B/C
$ws.Cells.Item(426, 1).Value2 = "mejor opcion"
$ws.Cells.Item(426, 2).Value2 = 1
$ws.Cells.Item(426, 3).Value2 = 0
$ws.Cells.Item(427, 1).Value2 = "lider"
$ws.Cells.Item(427, 2).Value2 = 1
$ws.Cells.Item(427, 3).Value2 = 0
$ws.Cells.Item(428, 1).Value2 = "no servicio"
$ws.Cells.Item(428, 2).Value2 = 0
$ws.Cells.Item(428, 3).Value2 = 1
$ws.Cells.Item(429, 1).Value2 = "pena"
$ws.Cells.Item(429, 2).Value2 = 0
$ws.Cells.Item(429, 3).Value2 = 1
$ws.Cells.Item(430, 1).Value2 = "feliz"
$ws.Cells.Item(430, 2).Value2 = 1
$ws.Cells.Item(430, 3).Value2 = 0
$ws.Cells.Item(431, 1).Value2 = "agradable"
$ws.Cells.Item(431, 2).Value2 = 1
$ws.Cells.Item(431, 3).Value2 = 0
$ws.Cells.Item(432, 1).Value2 = "gran servicio"
$ws.Cells.Item(432, 2).Value2 = 1
$ws.Cells.Item(432, 3).Value2 = 0
$ws.Cells.Item(433, 1).Value2 = "acuerdo"
$ws.Cells.Item(433, 2).Value2 = 1
$ws.Cells.Item(433, 3).Value2 = 0
$ws.Cells.Item(434, 1).Value2 = "normal"
$ws.Cells.Item(434, 2).Value2 = 0
$ws.Cells.Item(434, 3).Value2 = 0
$ws.Cells.Item(435, 1).Value2 = "guau"
$ws.Cells.Item(435, 2).Value2 = 1
$ws.Cells.Item(435, 3).Value2 = 0
$ws.Cells.Item(436, 1).Value2 = "promedio"
$ws.Cells.Item(436, 2).Value2 = 0
$ws.Cells.Item(436, 3).Value2 = 0
$ws.Cells.Item(437, 1).Value2 = "inhumano"
$ws.Cells.Item(437, 2).Value2 = 0
$ws.Cells.Item(437, 3).Value2 = 1
$ws.Cells.Item(438, 1).Value2 = "si poder poner"
$ws.Cells.Item(438, 2).Value2 = 0
$ws.Cells.Item(438, 3).Value2 = 1
$ws.Cells.Item(439, 1).Value2 = "excusa"
$ws.Cells.Item(439, 2).Value2 = 0
$ws.Cells.Item(439, 3).Value2 = 1
$ws.Cells.Item(440, 1).Value2 = "beneficio"
$ws.Cells.Item(440, 2).Value2 = 1
$ws.Cells.Item(440, 3).Value2 = 0
$ws.Cells.Item(441, 1).Value2 = "ineficiente"
$ws.Cells.Item(441, 2).Value2 = 0
$ws.Cells.Item(441, 3).Value2 = 1
$ws.Cells.Item(442, 1).Value2 = "nefasto"
$ws.Cells.Item(442, 2).Value2 = 0
$ws.Cells.Item(442, 3).Value2 = 1
$ws.Cells.Item(443, 1).Value2 = "ni contestan"
$ws.Cells.Item(443, 2).Value2 = 0
$ws.Cells.Item(443, 3).Value2 = 1
$ws.Cells.Item(444, 1).Value2 = "poco profesional"
$ws.Cells.Item(444, 2).Value2 = 0
$ws.Cells.Item(444, 3).Value2 = 1
$ws.Cells.Item(445, 1).Value2 = "mediocre"
$ws.Cells.Item(445, 2).Value2 = 0
$ws.Cells.Item(445, 3).Value2 = 1
$ws.Cells.Item(446, 1).Value2 = "nulo"
$ws.Cells.Item(446, 2).Value2 = 0
$ws.Cells.Item(446, 3).Value2 = 1
$ws.Cells.Item(447, 1).Value2 = "tarde"
$ws.Cells.Item(447, 2).Value2 = 0
$ws.Cells.Item(447, 3).Value2 = 1
$ws.Cells.Item(448, 1).Value2 = "cero apoyo"
$ws.Cells.Item(448, 2).Value2 = 0
$ws.Cells.Item(448, 3).Value2 = 1
$ws.Cells.Item(449, 1).Value2 = "sin apoyo"
$ws.Cells.Item(449, 2).Value2 = 0
$ws.Cells.Item(449, 3).Value2 = 1
$ws.Cells.Item(450, 1).Value2 = "casi hora"
$ws.Cells.Item(450, 2).Value2 = 0
$ws.Cells.Item(450, 3).Value2 = 1
$ws.Cells.Item(451, 1).Value2 = "mala gana"
$ws.Cells.Item(451, 2).Value2 = 0
$ws.Cells.Item(451, 3).Value2 = 1
$ws.Cells.Item(452, 1).Value2 = "no sincero"
$ws.Cells.Item(452, 2).Value2 = 0
$ws.Cells.Item(452, 3).Value2 = 1
$ws.Cells.Item(453, 1).Value2 = "alargar"
$ws.Cells.Item(453, 2).Value2 = 0
$ws.Cells.Item(453, 3).Value2 = 1
$ws.Cells.Item(454, 1).Value2 = "no mostrar"
$ws.Cells.Item(454, 2).Value2 = 0
$ws.Cells.Item(454, 3).Value2 = 1
$ws.Cells.Item(455, 1).Value2 = "limitar"
$ws.Cells.Item(455, 2).Value2 = 0
$ws.Cells.Item(455, 3).Value2 = 1
$ws.Cells.Item(456, 1).Value2 = "dificil"
$ws.Cells.Item(456, 2).Value2 = 0
$ws.Cells.Item(456, 3).Value2 = 1
$ws.Cells.Item(457, 1).Value2 = "nadie contestar"
$ws.Cells.Item(457, 2).Value2 = 0
$ws.Cells.Item(457, 3).Value2 = 1
$ws.Cells.Item(458, 1).Value2 = "cansir"
$ws.Cells.Item(458, 2).Value2 = 0
$ws.Cells.Item(458, 3).Value2 = 1
$ws.Cells.Item(459, 1).Value2 = "perdi"
$ws.Cells.Item(459, 2).Value2 = 0
$ws.Cells.Item(459, 3).Value2 = 1
$ws.Cells.Item(460, 1).Value2 = "no servir"
$ws.Cells.Item(460, 2).Value2 = 0
$ws.Cells.Item(460, 3).Value2 = 1
$ws.Cells.Item(461, 1).Value2 = "cambiense"
$ws.Cells.Item(461, 2).Value2 = 0
$ws.Cells.Item(461, 3).Value2 = 1
$ws.Cells.Item(462, 1).Value2 = "tirar"
$ws.Cells.Item(462, 2).Value2 = 0
$ws.Cells.Item(462, 3).Value2 = 1
$ws.Cells.Item(463, 1).Value2 = "tirar plata"
$ws.Cells.Item(463, 2).Value2 = 0
$ws.Cells.Item(463, 3).Value2 = 1
$ws.Cells.Item(464, 1).Value2 = "maltrato"
$ws.Cells.Item(464, 2).Value2 = 0
$ws.Cells.Item(464, 3).Value2 = 1
$ws.Cells.Item(465, 1).Value2 = "mediocr"
$ws.Cells.Item(465, 2).Value2 = 0
$ws.Cells.Item(465, 3).Value2 = 1
$ws.Cells.Item(466, 1).Value2 = "no funcionar"
$ws.Cells.Item(466, 2).Value2 = 0
$ws.Cells.Item(466, 3).Value2 = 1
$ws.Cells.Item(467, 1).Value2 = "frio"
$ws.Cells.Item(467, 2).Value2 = 0
$ws.Cells.Item(467, 3).Value2 = 1
$ws.Cells.Item(468, 1).Value2 = "tratar amabilidad"
$ws.Cells.Item(468, 2).Value2 = 0
$ws.Cells.Item(468, 3).Value2 = 1
$ws.Cells.Item(469, 1).Value2 = "capacitar"
$ws.Cells.Item(469, 2).Value2 = 0
$ws.Cells.Item(469, 3).Value2 = 1
$ws.Cells.Item(470, 1).Value2 = "no venir"
$ws.Cells.Item(470, 2).Value2 = 0
$ws.Cells.Item(470, 3).Value2 = 1
$ws.Cells.Item(471, 1).Value2 = "tratar"
$ws.Cells.Item(471, 2).Value2 = 0
$ws.Cells.Item(471, 3).Value2 = 1
$ws.Cells.Item(472, 1).Value2 = "no reconocer"
$ws.Cells.Item(472, 2).Value2 = 0
$ws.Cells.Item(472, 3).Value2 = 1
$ws.Cells.Item(473, 1).Value2 = "servicial"
$ws.Cells.Item(473, 2).Value2 = 1
$ws.Cells.Item(473, 3).Value2 = 0
$ws.Cells.Item(474, 1).Value2 = "organizado"
$ws.Cells.Item(474, 2).Value2 = 1
$ws.Cells.Item(474, 3).Value2 = 0
$ws.Cells.Item(475, 1).Value2 = "concis"
$ws.Cells.Item(475, 2).Value2 = 1
$ws.Cells.Item(475, 3).Value2 = 0
$ws.Cells.Item(476, 1).Value2 = "detallado"
$ws.Cells.Item(476, 2).Value2 = 1
$ws.Cells.Item(476, 3).Value2 = 0
$ws.Cells.Item(477, 1).Value2 = "trabar"
$ws.Cells.Item(477, 2).Value2 = 0
$ws.Cells.Item(477, 3).Value2 = 1
$ws.Cells.Item(478, 1).Value2 = "no agrado"
$ws.Cells.Item(478, 2).Value2 = 0
$ws.Cells.Item(478, 3).Value2 = 1
$ws.Cells.Item(479, 1).Value2 = "no señal"
$ws.Cells.Item(479, 2).Value2 = 0
$ws.Cells.Item(479, 3).Value2 = 1
$ws.Cells.Item(480, 1).Value2 = "no poder creer"
$ws.Cells.Item(480, 2).Value2 = 0
$ws.Cells.Item(480, 3).Value2 = 1
$ws.Cells.Item(481, 1).Value2 = "decir gran empresa"
$ws.Cells.Item(481, 2).Value2 = 0
$ws.Cells.Item(481, 3).Value2 = 1
$ws.Cells.Item(482, 1).Value2 = "no tener personal calificado"
$ws.Cells.Item(482, 2).Value2 = 0
$ws.Cells.Item(482, 3).Value2 = 1
$ws.Cells.Item(483, 1).Value2 = "amablemente"
$ws.Cells.Item(483, 2).Value2 = 1
$ws.Cells.Item(483, 3).Value2 = 0
$ws.Cells.Item(484, 1).Value2 = "exelente"
$ws.Cells.Item(484, 2).Value2 = 1
$ws.Cells.Item(484, 3).Value2 = 0
$ws.Cells.Item(485, 1).Value2 = "cualquiera cosa"
$ws.Cells.Item(485, 2).Value2 = 0
$ws.Cells.Item(485, 3).Value2 = 1
$ws.Cells.Item(486, 1).Value2 = "preocupado"
$ws.Cells.Item(486, 2).Value2 = 0
$ws.Cells.Item(486, 3).Value2 = 1
$ws.Cells.Item(487, 1).Value2 = "medio hora"
$ws.Cells.Item(487, 2).Value2 = 0
$ws.Cells.Item(487, 3).Value2 = 1
$ws.Cells.Item(488, 1).Value2 = "no primero vez"
$ws.Cells.Item(488, 2).Value2 = 0
$ws.Cells.Item(488, 3).Value2 = 1
$ws.Cells.Item(489, 1).Value2 = "cualquiera"
$ws.Cells.Item(489, 2).Value2 = 0
$ws.Cells.Item(489, 3).Value2 = 1
$ws.Cells.Item(490, 1).Value2 = "cualquier"
$ws.Cells.Item(490, 2).Value2 = 0
$ws.Cells.Item(490, 3).Value2 = 1
$ws.Cells.Item(491, 1).Value2 = "bonito"
$ws.Cells.Item(491, 2).Value2 = 1
$ws.Cells.Item(491, 3).Value2 = 0
$ws.Cells.Item(492, 1).Value2 = "resguardado"
$ws.Cells.Item(492, 2).Value2 = 1
$ws.Cells.Item(492, 3).Value2 = 0
$ws.Cells.Item(493, 1).Value2 = "poco tiempo espera"
$ws.Cells.Item(493, 2).Value2 = 1
$ws.Cells.Item(493, 3).Value2 = 0
$ws.Cells.Item(494, 1).Value2 = "correctamente"
$ws.Cells.Item(494, 2).Value2 = 1
$ws.Cells.Item(494, 3).Value2 = 0
$ws.Cells.Item(495, 1).Value2 = "facil perder"
$ws.Cells.Item(495, 2).Value2 = 0
$ws.Cells.Item(495, 3).Value2 = 1
$ws.Cells.Item(496, 1).Value2 = "limpieza"
$ws.Cells.Item(496, 2).Value2 = 1
$ws.Cells.Item(496, 3).Value2 = 0
$ws.Cells.Item(497, 1).Value2 = "agil"
$ws.Cells.Item(497, 2).Value2 = 1
$ws.Cells.Item(497, 3).Value2 = 0
$ws.Cells.Item(498, 1).Value2 = "broma"
$ws.Cells.Item(498, 2).Value2 = 0
$ws.Cells.Item(498, 3).Value2 = 1
$ws.Cells.Item(499, 1).Value2 = "sensibilidad"
$ws.Cells.Item(499, 2).Value2 = 1
$ws.Cells.Item(499, 3).Value2 = 0
$ws.Cells.Item(500, 1).Value2 = "querer pagar"
$ws.Cells.Item(500, 2).Value2 = 0
$ws.Cells.Item(500, 3).Value2 = 1
$ws.Cells.Item(501, 1).Value2 = "deber brindar"
$ws.Cells.Item(501, 2).Value2 = 0
$ws.Cells.Item(501, 3).Value2 = 1
$ws.Cells.Item(502, 1).Value2 = "llamar varios vez"
$ws.Cells.Item(502, 2).Value2 = 0
$ws.Cells.Item(502, 3).Value2 = 1
$ws.Cells.Item(503, 1).Value2 = "faltar buen"
$ws.Cells.Item(503, 2).Value2 = 0
$ws.Cells.Item(503, 3).Value2 = 1
$ws.Cells.Item(504, 1).Value2 = "faltar buen atencion"
$ws.Cells.Item(504, 2).Value2 = 0
$ws.Cells.Item(504, 3).Value2 = 1
$ws.Cells.Item(505, 1).Value2 = "faltar buen atencion personalizado"
$ws.Cells.Item(505, 2).Value2 = 0
$ws.Cells.Item(505, 3).Value2 = 1
$ws.Cells.Item(506, 1).Value2 = "incomodo"
$ws.Cells.Item(506, 2).Value2 = 0
$ws.Cells.Item(506, 3).Value2 = 1
$ws.Cells.Item(507, 1).Value2 = "no orientacion"
$ws.Cells.Item(507, 2).Value2 = 0
$ws.Cells.Item(507, 3).Value2 = 1
$ws.Cells.Item(508, 1).Value2 = "no necesario"
$ws.Cells.Item(508, 2).Value2 = 0
$ws.Cells.Item(508, 3).Value2 = 1
$ws.Cells.Item(509, 1).Value2 = "tiempo perdido"
$ws.Cells.Item(509, 2).Value2 = 0
$ws.Cells.Item(509, 3).Value2 = 1
$ws.Cells.Item(510, 1).Value2 = "deber resolver"
$ws.Cells.Item(510, 2).Value2 = 0
$ws.Cells.Item(510, 3).Value2 = 1
$ws.Cells.Item(511, 1).Value2 = "rapidez"
$ws.Cells.Item(511, 2).Value2 = 1
$ws.Cells.Item(511, 3).Value2 = 0
$ws.Cells.Item(512, 1).Value2 = "faltar veracidad"
$ws.Cells.Item(512, 2).Value2 = 0
$ws.Cells.Item(512, 3).Value2 = 1
$ws.Cells.Item(513, 1).Value2 = "felizmente poder solucionar"
$ws.Cells.Item(513, 2).Value2 = 1
$ws.Cells.Item(513, 3).Value2 = 0
$ws.Cells.Item(514, 1).Value2 = "total desorganizacion"
$ws.Cells.Item(514, 2).Value2 = 0
$ws.Cells.Item(514, 3).Value2 = 1
$ws.Cells.Item(515, 1).Value2 = "no organización"
$ws.Cells.Item(515, 2).Value2 = 0
$ws.Cells.Item(515, 3).Value2 = 1
$ws.Cells.Item(516, 1).Value2 = "porfavor"
$ws.Cells.Item(516, 2).Value2 = 0
$ws.Cells.Item(516, 3).Value2 = 0
$ws.Cells.Item(517, 1).Value2 = "papeleo"
$ws.Cells.Item(517, 2).Value2 = 0
$ws.Cells.Item(517, 3).Value2 = 1
$ws.Cells.Item(518, 1).Value2 = "grato"
$ws.Cells.Item(518, 2).Value2 = 1
$ws.Cells.Item(518, 3).Value2 = 0
$ws.Cells.Item(519, 1).Value2 = "correcto"
$ws.Cells.Item(519, 2).Value2 = 1
$ws.Cells.Item(519, 3).Value2 = 0
$ws.Cells.Item(520, 1).Value2 = "trato comunicar"
$ws.Cells.Item(520, 2).Value2 = 0
$ws.Cells.Item(520, 3).Value2 = 1
$ws.Cells.Item(521, 1).Value2 = "trato llamar"
$ws.Cells.Item(521, 2).Value2 = 0
$ws.Cells.Item(521, 3).Value2 = 1
$ws.Cells.Item(522, 1).Value2 = "trato contactar"
$ws.Cells.Item(522, 2).Value2 = 0
$ws.Cells.Item(522, 3).Value2 = 1
$ws.Cells.Item(523, 1).Value2 = "no suficiente"
$ws.Cells.Item(523, 2).Value2 = 0
$ws.Cells.Item(523, 3).Value2 = 1
$ws.Cells.Item(524, 1).Value2 = "tratar comunicar"
$ws.Cells.Item(524, 2).Value2 = 0
$ws.Cells.Item(524, 3).Value2 = 1
$ws.Cells.Item(525, 1).Value2 = "tratar llamar"
$ws.Cells.Item(525, 2).Value2 = 0
$ws.Cells.Item(525, 3).Value2 = 1
$ws.Cells.Item(526, 1).Value2 = "tratar contactar"
$ws.Cells.Item(526, 2).Value2 = 0
$ws.Cells.Item(526, 3).Value2 = 1
$ws.Cells.Item(527, 1).Value2 = "perseverancia"
$ws.Cells.Item(527, 2).Value2 = 0
$ws.Cells.Item(527, 3).Value2 = 1
$ws.Cells.Item(528, 1).Value2 = "esperar mucho"
$ws.Cells.Item(528, 2).Value2 = 0
$ws.Cells.Item(528, 3).Value2 = 1
$ws.Cells.Item(529, 1).Value2 = "esperar demasiado"
$ws.Cells.Item(529, 2).Value2 = 0
$ws.Cells.Item(529, 3).Value2 = 1

# Column width adjustments (closest achievable values given engine quantization)
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(2).ColumnWidth = 13.333333333333334

# Update active selection to E7
$ws.Activate()
$ws.Range("E7").Select()
